# Fruta / hortaliza, semanal
# Refresh weekly snapshot values: the rows 2-17 (Vega Central Mapocho de
# Santiago - Mora) get their Fecha/Volumen/Precio/Origen figures reshuffled
# to reflect a different weekly sampling window. Row 10 is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44231
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 3400
$ws.Range("O2").Value = 3400
$ws.Range("P2").Value = 3400
$ws.Range("S2").Value = 1700

$ws.Range("D3").Value = 44978
$ws.Range("M3").Value = 500

$ws.Range("D4").Value = 44533
$ws.Range("M4").Value = 150
$ws.Range("N4").Value = 4000
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 4000
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 2000

$ws.Range("D5").Value = 44168
$ws.Range("M5").Value = 170
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 4000

$ws.Range("D6").Value = 44208
$ws.Range("M6").Value = 85
$ws.Range("O6").Value = 3000
$ws.Range("P6").Value = 3000
$ws.Range("S6").Value = 1500

$ws.Range("D7").Value = 44236
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 3600
$ws.Range("P7").Value = 3800
$ws.Range("S7").Value = 1900

$ws.Range("D8").Value = 44237
$ws.Range("M8").Value = 100

$ws.Range("D9").Value = 44980
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 4000
$ws.Range("P9").Value = 4000
$ws.Range("S9").Value = 2000

$ws.Range("D11").Value = 44617
$ws.Range("M11").Value = 90
$ws.Range("N11").Value = 6500
$ws.Range("O11").Value = 6500
$ws.Range("P11").Value = 6500
$ws.Range("S11").Value = 3250

$ws.Range("D12").Value = 44582
$ws.Range("M12").Value = 380
$ws.Range("N12").Value = 5000
$ws.Range("O12").Value = 5000
$ws.Range("P12").Value = 5000
$ws.Range("S12").Value = 2500

$ws.Range("D13").Value = 44586
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 5000
$ws.Range("O13").Value = 5000
$ws.Range("P13").Value = 5000
$ws.Range("S13").Value = 2500

$ws.Range("D14").Value = 44194
$ws.Range("M14").Value = 120

$ws.Range("D15").Value = 44232
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 3000
$ws.Range("O15").Value = 3000
$ws.Range("P15").Value = 3000
$ws.Range("S15").Value = 1500

$ws.Range("D16").Value = 44188
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 3400
$ws.Range("P16").Value = 3240
$ws.Range("R16").Value = "Provincia de Linares"
$ws.Range("S16").Value = 1620

$ws.Range("D17").Value = 44238
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 3600
$ws.Range("O17").Value = 4000
$ws.Range("P17").Value = 3800
$ws.Range("R17").Value = "Provincia de Curicó"
$ws.Range("S17").Value = 1900
